{"js": "// Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings in\n// the practice-table cells with their updated values, in document order.\n// (Every populated data cell in the table changed \u2014 one text replacement\n// per cell, old text -> new text.)\nconst replacements = [\n  [\"529\u00f76=88, 1\", \"342\u00f77=48, 6\"],\n  [\"852\u00f77=121, 5\", \"849\u00f77=121, 2\"],\n  [\"372\u00f79=41, 3\", \"682\u00f73=227, 1\"],\n  [\"562\u00f72=281, 0\", \"367\u00f72=183, 1\"],\n  [\"635\u00f73=211, 2\", \"576\u00f77=82, 2\"],\n  [\"871\u00f77=124, 3\", \"570\u00f74=142, 2\"],\n  [\"476\u00f79=52, 8\", \"917\u00f75=183, 2\"],\n  [\"901\u00f72=450, 1\", \"575\u00f75=115, 0\"],\n  [\"220\u00f75=44, 0\", \"188\u00f79=20, 8\"],\n  [\"863\u00f77=123, 2\", \"466\u00f72=233, 0\"],\n  [\"173\u00f73=57, 2\", \"847\u00f77=121, 0\"],\n  [\"844\u00f75=168, 4\", \"507\u00f77=72, 3\"],\n  [\"187\u00f73=62, 1\", \"920\u00f72=460, 0\"],\n  [\"960\u00f78=120, 0\", \"421\u00f73=140, 1\"],\n  [\"525\u00f74=131, 1\", \"804\u00f72=402, 0\"],\n  [\"934\u00f76=155, 4\", \"888\u00f74=222, 0\"],\n  [\"761\u00f72=380, 1\", \"711\u00f78=88, 7\"],\n  [\"383\u00f73=127, 2\", \"396\u00f78=49, 4\"],\n  [\"802\u00f78=100, 2\", \"115\u00f78=14, 3\"],\n  [\"650\u00f79=72, 2\", \"538\u00f72=269, 0\"],\n  [\"866\u00f73=288, 2\", \"784\u00f79=87, 1\"],\n  [\"547\u00f74=136, 3\", \"608\u00f76=101, 2\"],\n  [\"347\u00f76=57, 5\", \"399\u00f75=79, 4\"],\n  [\"659\u00f72=329, 1\", \"442\u00f75=88, 2\"],\n  [\"376\u00f75=75, 1\", \"431\u00f72=215, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"dividend\u00f7divisor=quotient, remainder\" answer strings in\n# the practice-table cells with their updated values, in document order.\n# (Every populated data cell in the table changed -- one Find/Replace per\n# cell, old text -> new text.)\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"529\u00f76=88, 1\", \"342\u00f77=48, 6\"),\n    @(\"852\u00f77=121, 5\", \"849\u00f77=121, 2\"),\n    @(\"372\u00f79=41, 3\", \"682\u00f73=227, 1\"),\n    @(\"562\u00f72=281, 0\", \"367\u00f72=183, 1\"),\n    @(\"635\u00f73=211, 2\", \"576\u00f77=82, 2\"),\n    @(\"871\u00f77=124, 3\", \"570\u00f74=142, 2\"),\n    @(\"476\u00f79=52, 8\", \"917\u00f75=183, 2\"),\n    @(\"901\u00f72=450, 1\", \"575\u00f75=115, 0\"),\n    @(\"220\u00f75=44, 0\", \"188\u00f79=20, 8\"),\n    @(\"863\u00f77=123, 2\", \"466\u00f72=233, 0\"),\n    @(\"173\u00f73=57, 2\", \"847\u00f77=121, 0\"),\n    @(\"844\u00f75=168, 4\", \"507\u00f77=72, 3\"),\n    @(\"187\u00f73=62, 1\", \"920\u00f72=460, 0\"),\n    @(\"960\u00f78=120, 0\", \"421\u00f73=140, 1\"),\n    @(\"525\u00f74=131, 1\", \"804\u00f72=402, 0\"),\n    @(\"934\u00f76=155, 4\", \"888\u00f74=222, 0\"),\n    @(\"761\u00f72=380, 1\", \"711\u00f78=88, 7\"),\n    @(\"383\u00f73=127, 2\", \"396\u00f78=49, 4\"),\n    @(\"802\u00f78=100, 2\", \"115\u00f78=14, 3\"),\n    @(\"650\u00f79=72, 2\", \"538\u00f72=269, 0\"),\n    @(\"866\u00f73=288, 2\", \"784\u00f79=87, 1\"),\n    @(\"547\u00f74=136, 3\", \"608\u00f76=101, 2\"),\n    @(\"347\u00f76=57, 5\", \"399\u00f75=79, 4\"),\n    @(\"659\u00f72=329, 1\", \"442\u00f75=88, 2\"),\n    @(\"376\u00f75=75, 1\", \"431\u00f72=215, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
